$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit shuffles the per-row data (Fecha, Variedad, Calidad, Volumen, prices)
# across rows 2-18. Capture the original values first (columns D,H,I,J,K,L,M,P),
# then write each target row from its mapped source row, so the operation is
# correct regardless of write order.
$orig = @{}
$orig[2] = @{
    D = $ws.Range("D2").Value2
    H = $ws.Range("H2").Value2
    I = $ws.Range("I2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    P = $ws.Range("P2").Value2
}
$orig[3] = @{
    D = $ws.Range("D3").Value2
    H = $ws.Range("H3").Value2
    I = $ws.Range("I3").Value2
    J = $ws.Range("J3").Value2
    K = $ws.Range("K3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    P = $ws.Range("P3").Value2
}
$orig[4] = @{
    D = $ws.Range("D4").Value2
    H = $ws.Range("H4").Value2
    I = $ws.Range("I4").Value2
    J = $ws.Range("J4").Value2
    K = $ws.Range("K4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    P = $ws.Range("P4").Value2
}
$orig[5] = @{
    D = $ws.Range("D5").Value2
    H = $ws.Range("H5").Value2
    I = $ws.Range("I5").Value2
    J = $ws.Range("J5").Value2
    K = $ws.Range("K5").Value2
    L = $ws.Range("L5").Value2
    M = $ws.Range("M5").Value2
    P = $ws.Range("P5").Value2
}
$orig[6] = @{
    D = $ws.Range("D6").Value2
    H = $ws.Range("H6").Value2
    I = $ws.Range("I6").Value2
    J = $ws.Range("J6").Value2
    K = $ws.Range("K6").Value2
    L = $ws.Range("L6").Value2
    M = $ws.Range("M6").Value2
    P = $ws.Range("P6").Value2
}
$orig[7] = @{
    D = $ws.Range("D7").Value2
    H = $ws.Range("H7").Value2
    I = $ws.Range("I7").Value2
    J = $ws.Range("J7").Value2
    K = $ws.Range("K7").Value2
    L = $ws.Range("L7").Value2
    M = $ws.Range("M7").Value2
    P = $ws.Range("P7").Value2
}
$orig[8] = @{
    D = $ws.Range("D8").Value2
    H = $ws.Range("H8").Value2
    I = $ws.Range("I8").Value2
    J = $ws.Range("J8").Value2
    K = $ws.Range("K8").Value2
    L = $ws.Range("L8").Value2
    M = $ws.Range("M8").Value2
    P = $ws.Range("P8").Value2
}
$orig[9] = @{
    D = $ws.Range("D9").Value2
    H = $ws.Range("H9").Value2
    I = $ws.Range("I9").Value2
    J = $ws.Range("J9").Value2
    K = $ws.Range("K9").Value2
    L = $ws.Range("L9").Value2
    M = $ws.Range("M9").Value2
    P = $ws.Range("P9").Value2
}
$orig[10] = @{
    D = $ws.Range("D10").Value2
    H = $ws.Range("H10").Value2
    I = $ws.Range("I10").Value2
    J = $ws.Range("J10").Value2
    K = $ws.Range("K10").Value2
    L = $ws.Range("L10").Value2
    M = $ws.Range("M10").Value2
    P = $ws.Range("P10").Value2
}
$orig[11] = @{
    D = $ws.Range("D11").Value2
    H = $ws.Range("H11").Value2
    I = $ws.Range("I11").Value2
    J = $ws.Range("J11").Value2
    K = $ws.Range("K11").Value2
    L = $ws.Range("L11").Value2
    M = $ws.Range("M11").Value2
    P = $ws.Range("P11").Value2
}
$orig[12] = @{
    D = $ws.Range("D12").Value2
    H = $ws.Range("H12").Value2
    I = $ws.Range("I12").Value2
    J = $ws.Range("J12").Value2
    K = $ws.Range("K12").Value2
    L = $ws.Range("L12").Value2
    M = $ws.Range("M12").Value2
    P = $ws.Range("P12").Value2
}
$orig[13] = @{
    D = $ws.Range("D13").Value2
    H = $ws.Range("H13").Value2
    I = $ws.Range("I13").Value2
    J = $ws.Range("J13").Value2
    K = $ws.Range("K13").Value2
    L = $ws.Range("L13").Value2
    M = $ws.Range("M13").Value2
    P = $ws.Range("P13").Value2
}
$orig[14] = @{
    D = $ws.Range("D14").Value2
    H = $ws.Range("H14").Value2
    I = $ws.Range("I14").Value2
    J = $ws.Range("J14").Value2
    K = $ws.Range("K14").Value2
    L = $ws.Range("L14").Value2
    M = $ws.Range("M14").Value2
    P = $ws.Range("P14").Value2
}
$orig[15] = @{
    D = $ws.Range("D15").Value2
    H = $ws.Range("H15").Value2
    I = $ws.Range("I15").Value2
    J = $ws.Range("J15").Value2
    K = $ws.Range("K15").Value2
    L = $ws.Range("L15").Value2
    M = $ws.Range("M15").Value2
    P = $ws.Range("P15").Value2
}
$orig[16] = @{
    D = $ws.Range("D16").Value2
    H = $ws.Range("H16").Value2
    I = $ws.Range("I16").Value2
    J = $ws.Range("J16").Value2
    K = $ws.Range("K16").Value2
    L = $ws.Range("L16").Value2
    M = $ws.Range("M16").Value2
    P = $ws.Range("P16").Value2
}
$orig[17] = @{
    D = $ws.Range("D17").Value2
    H = $ws.Range("H17").Value2
    I = $ws.Range("I17").Value2
    J = $ws.Range("J17").Value2
    K = $ws.Range("K17").Value2
    L = $ws.Range("L17").Value2
    M = $ws.Range("M17").Value2
    P = $ws.Range("P17").Value2
}
$orig[18] = @{
    D = $ws.Range("D18").Value2
    H = $ws.Range("H18").Value2
    I = $ws.Range("I18").Value2
    J = $ws.Range("J18").Value2
    K = $ws.Range("K18").Value2
    L = $ws.Range("L18").Value2
    M = $ws.Range("M18").Value2
    P = $ws.Range("P18").Value2
}

# Apply the shuffled values: each target row gets the original values of its mapped source row.
$ws.Range("D2").Value2 = $orig[10].D
$ws.Range("H2").Value2 = $orig[10].H
$ws.Range("I2").Value2 = $orig[10].I
$ws.Range("J2").Value2 = $orig[10].J
$ws.Range("K2").Value2 = $orig[10].K
$ws.Range("L2").Value2 = $orig[10].L
$ws.Range("M2").Value2 = $orig[10].M
$ws.Range("P2").Value2 = $orig[10].P

$ws.Range("D3").Value2 = $orig[12].D
$ws.Range("H3").Value2 = $orig[12].H
$ws.Range("I3").Value2 = $orig[12].I
$ws.Range("J3").Value2 = $orig[12].J
$ws.Range("K3").Value2 = $orig[12].K
$ws.Range("L3").Value2 = $orig[12].L
$ws.Range("M3").Value2 = $orig[12].M
$ws.Range("P3").Value2 = $orig[12].P

$ws.Range("D4").Value2 = $orig[18].D
$ws.Range("H4").Value2 = $orig[18].H
$ws.Range("I4").Value2 = $orig[18].I
$ws.Range("J4").Value2 = $orig[18].J
$ws.Range("K4").Value2 = $orig[18].K
$ws.Range("L4").Value2 = $orig[18].L
$ws.Range("M4").Value2 = $orig[18].M
$ws.Range("P4").Value2 = $orig[18].P

$ws.Range("D5").Value2 = $orig[2].D
$ws.Range("H5").Value2 = $orig[2].H
$ws.Range("I5").Value2 = $orig[2].I
$ws.Range("J5").Value2 = $orig[2].J
$ws.Range("K5").Value2 = $orig[2].K
$ws.Range("L5").Value2 = $orig[2].L
$ws.Range("M5").Value2 = $orig[2].M
$ws.Range("P5").Value2 = $orig[2].P

$ws.Range("D6").Value2 = $orig[8].D
$ws.Range("H6").Value2 = $orig[8].H
$ws.Range("I6").Value2 = $orig[8].I
$ws.Range("J6").Value2 = $orig[8].J
$ws.Range("K6").Value2 = $orig[8].K
$ws.Range("L6").Value2 = $orig[8].L
$ws.Range("M6").Value2 = $orig[8].M
$ws.Range("P6").Value2 = $orig[8].P

$ws.Range("D7").Value2 = $orig[17].D
$ws.Range("H7").Value2 = $orig[17].H
$ws.Range("I7").Value2 = $orig[17].I
$ws.Range("J7").Value2 = $orig[17].J
$ws.Range("K7").Value2 = $orig[17].K
$ws.Range("L7").Value2 = $orig[17].L
$ws.Range("M7").Value2 = $orig[17].M
$ws.Range("P7").Value2 = $orig[17].P

$ws.Range("D8").Value2 = $orig[7].D
$ws.Range("H8").Value2 = $orig[7].H
$ws.Range("I8").Value2 = $orig[7].I
$ws.Range("J8").Value2 = $orig[7].J
$ws.Range("K8").Value2 = $orig[7].K
$ws.Range("L8").Value2 = $orig[7].L
$ws.Range("M8").Value2 = $orig[7].M
$ws.Range("P8").Value2 = $orig[7].P

$ws.Range("D9").Value2 = $orig[4].D
$ws.Range("H9").Value2 = $orig[4].H
$ws.Range("I9").Value2 = $orig[4].I
$ws.Range("J9").Value2 = $orig[4].J
$ws.Range("K9").Value2 = $orig[4].K
$ws.Range("L9").Value2 = $orig[4].L
$ws.Range("M9").Value2 = $orig[4].M
$ws.Range("P9").Value2 = $orig[4].P

$ws.Range("D10").Value2 = $orig[15].D
$ws.Range("H10").Value2 = $orig[15].H
$ws.Range("I10").Value2 = $orig[15].I
$ws.Range("J10").Value2 = $orig[15].J
$ws.Range("K10").Value2 = $orig[15].K
$ws.Range("L10").Value2 = $orig[15].L
$ws.Range("M10").Value2 = $orig[15].M
$ws.Range("P10").Value2 = $orig[15].P

$ws.Range("D11").Value2 = $orig[9].D
$ws.Range("H11").Value2 = $orig[9].H
$ws.Range("I11").Value2 = $orig[9].I
$ws.Range("J11").Value2 = $orig[9].J
$ws.Range("K11").Value2 = $orig[9].K
$ws.Range("L11").Value2 = $orig[9].L
$ws.Range("M11").Value2 = $orig[9].M
$ws.Range("P11").Value2 = $orig[9].P

$ws.Range("D12").Value2 = $orig[13].D
$ws.Range("H12").Value2 = $orig[13].H
$ws.Range("I12").Value2 = $orig[13].I
$ws.Range("J12").Value2 = $orig[13].J
$ws.Range("K12").Value2 = $orig[13].K
$ws.Range("L12").Value2 = $orig[13].L
$ws.Range("M12").Value2 = $orig[13].M
$ws.Range("P12").Value2 = $orig[13].P

$ws.Range("D13").Value2 = $orig[3].D
$ws.Range("H13").Value2 = $orig[3].H
$ws.Range("I13").Value2 = $orig[3].I
$ws.Range("J13").Value2 = $orig[3].J
$ws.Range("K13").Value2 = $orig[3].K
$ws.Range("L13").Value2 = $orig[3].L
$ws.Range("M13").Value2 = $orig[3].M
$ws.Range("P13").Value2 = $orig[3].P

$ws.Range("D14").Value2 = $orig[5].D
$ws.Range("H14").Value2 = $orig[5].H
$ws.Range("I14").Value2 = $orig[5].I
$ws.Range("J14").Value2 = $orig[5].J
$ws.Range("K14").Value2 = $orig[5].K
$ws.Range("L14").Value2 = $orig[5].L
$ws.Range("M14").Value2 = $orig[5].M
$ws.Range("P14").Value2 = $orig[5].P

$ws.Range("D15").Value2 = $orig[11].D
$ws.Range("H15").Value2 = $orig[11].H
$ws.Range("I15").Value2 = $orig[11].I
$ws.Range("J15").Value2 = $orig[11].J
$ws.Range("K15").Value2 = $orig[11].K
$ws.Range("L15").Value2 = $orig[11].L
$ws.Range("M15").Value2 = $orig[11].M
$ws.Range("P15").Value2 = $orig[11].P

$ws.Range("D16").Value2 = $orig[6].D
$ws.Range("H16").Value2 = $orig[6].H
$ws.Range("I16").Value2 = $orig[6].I
$ws.Range("J16").Value2 = $orig[6].J
$ws.Range("K16").Value2 = $orig[6].K
$ws.Range("L16").Value2 = $orig[6].L
$ws.Range("M16").Value2 = $orig[6].M
$ws.Range("P16").Value2 = $orig[6].P

$ws.Range("D17").Value2 = $orig[16].D
$ws.Range("H17").Value2 = $orig[16].H
$ws.Range("I17").Value2 = $orig[16].I
$ws.Range("J17").Value2 = $orig[16].J
$ws.Range("K17").Value2 = $orig[16].K
$ws.Range("L17").Value2 = $orig[16].L
$ws.Range("M17").Value2 = $orig[16].M
$ws.Range("P17").Value2 = $orig[16].P

$ws.Range("D18").Value2 = $orig[14].D
$ws.Range("H18").Value2 = $orig[14].H
$ws.Range("I18").Value2 = $orig[14].I
$ws.Range("J18").Value2 = $orig[14].J
$ws.Range("K18").Value2 = $orig[14].K
$ws.Range("L18").Value2 = $orig[14].L
$ws.Range("M18").Value2 = $orig[14].M
$ws.Range("P18").Value2 = $orig[14].P

